{"js": "// Reto71.docx: insert \" separados por espacio\" after \"...riesgo MEDIO\" and\n// before \", en caso de no haber ninguno devolver NA.\" in the bullet list item\n// that asks to indicate the IRCA rating for MEDIO-risk water bodies.\n\nconst oldText =\n  \"Indicar la calificaci\u00f3n IRCA de los cuerpos de agua que tienen un nivel de riesgo MEDIO, en caso de no haber ninguno devolver NA.\";\nconst newText =\n  \"Indicar la calificaci\u00f3n IRCA de los cuerpos de agua que tienen un nivel de riesgo MEDIO separados por espacio, en caso de no haber ninguno devolver NA.\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target paragraph text not found: \" + oldText);\n}\n\n// Replace the whole paragraph run's text with the updated wording in place,\n// preserving the existing run formatting (color, etc.).\nresults.items[0].insertText(newText, \"Replace\");\nawait context.sync();\n", "ps1": "# Reto71.docx: insert \" separados por espacio\" after \"...riesgo MEDIO\" and\n# before \", en caso de no haber ninguno devolver NA.\" in the bullet list item\n# that asks to indicate the IRCA rating for MEDIO-risk water bodies.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"MEDIO, en caso de no haber ninguno devolver NA.\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"MEDIO separados por espacio, en caso de no haber ninguno devolver NA.\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
